$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D as Text for the data rows so that numeric-looking
# strings (e.g. "1.000", "0.9968") are preserved as text, matching the
# original inline-string storage instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '24.755.40'
$ws.Range("E2").Value = '  +2.15%  '

# Row 3
$ws.Range("D3").Value = '1.703.92'
$ws.Range("E3").Value = '  +1.22%  '

# Row 4
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.32%  '

# Row 5
$ws.Range("D5").Value = '309.29'
$ws.Range("E5").Value = '  -0.02%  '

# Row 6
$ws.Range("D6").Value = '0.9968'
$ws.Range("E6").Value = '  -0.16%  '

# Row 7
$ws.Range("D7").Value = '0.3742'
$ws.Range("E7").Value = '  +0.09%  '

# Row 8
$ws.Range("D8").Value = '49.16'
$ws.Range("E8").Value = '  +2.70%  '

# Row 9
$ws.Range("D9").Value = '0.3441'
$ws.Range("E9").Value = '  -0.54%  '

# Row 10
$ws.Range("D10").Value = '1.188'
$ws.Range("E10").Value = '  -0.39%  '

# Row 11
$ws.Range("D11").Value = '0.07463'
$ws.Range("E11").Value = '  +1.99%  '

# Row 12
$ws.Range("D12").Value = '0.9961'
$ws.Range("E12").Value = '  -0.41%  '

# Row 13
$ws.Range("D13").Value = '20.90'
$ws.Range("E13").Value = '  +1.84%  '

# Row 14
$ws.Range("D14").Value = '6.243'
$ws.Range("E14").Value = '  +1.66%  '

# Row 15
$ws.Range("D15").Value = '6.938'
$ws.Range("E15").Value = '  +2.03%  '

# Row 16
$ws.Range("D16").Value = '1.705.43'
$ws.Range("E16").Value = '  +1.23%  '

# Row 17
$ws.Range("D17").Value = '0.00001124'
$ws.Range("E17").Value = '  +0.98%  '

# Row 18
$ws.Range("D18").Value = '0.06698'
$ws.Range("E18").Value = '  -0.46%  '

# Row 19
$ws.Range("D19").Value = '0.9971'
$ws.Range("E19").Value = '  -0.12%  '

# Row 20
$ws.Range("D20").Value = '83.70'
$ws.Range("E20").Value = '  +1.70%  '

# Row 21
$ws.Range("D21").Value = '17.10'
$ws.Range("E21").Value = '  +3.41%  '

# Row 22
$ws.Range("D22").Value = '6.337'
$ws.Range("E22").Value = '  +3.53%  '

# Row 23
$ws.Range("D23").Value = '13.13'
$ws.Range("E23").Value = '  +8.85%  '

# Row 24
$ws.Range("D24").Value = '24.724.30'
$ws.Range("E24").Value = '  +2.12%  '

# Row 25
$ws.Range("D25").Value = '2.418'
$ws.Range("E25").Value = '  -0.05%  '

# Row 26
$ws.Range("D26").Value = '2.775'
$ws.Range("E26").Value = '  +3.30%  '

# Row 27
$ws.Range("D27").Value = '20.12'
$ws.Range("E27").Value = '  +2.19%  '

# Row 28
$ws.Range("D28").Value = '150.26'
$ws.Range("E28").Value = '  -2.37%  '

# Row 29
$ws.Range("D29").Value = '130.97'
$ws.Range("E29").Value = '  +2.87%  '

# Row 30
$ws.Range("D30").Value = '1.895.17'
$ws.Range("E30").Value = '  +1.38%  '

# Row 31
$ws.Range("D31").Value = '1.189'
$ws.Range("E31").Value = '  +19.92%  '

# Row 32
$ws.Range("D32").Value = '6.790'
$ws.Range("E32").Value = '  +4.99%  '

# Row 33
$ws.Range("D33").Value = '4.185'
$ws.Range("E33").Value = '  +1.95%  '

# Row 34
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '13.69'
$ws.Range("E34").Value = '  +9.47%  '

# Row 35
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '1.779'
$ws.Range("E35").Value = '  -0.44%  '

# Row 36
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").Value = '0.08795'
$ws.Range("E36").Value = '  +3.69%  '

# Row 37
$ws.Range("D37").Value = '5.530'
$ws.Range("E37").Value = '  +2.53%  '

# Row 38
$ws.Range("D38").Value = '0.06526'
$ws.Range("E38").Value = '  +0.27%  '

# Row 39
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '8.942'
$ws.Range("E39").Value = '  -0.44%  '

# Row 40
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.2234'
$ws.Range("E40").Value = '  +4.37%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.02382'
$ws.Range("E41").Value = '  +0.99%  '

# Row 42
$ws.Range("D42").Value = '1.275'
$ws.Range("E42").Value = '  -0.33%  '

# Row 43
$ws.Range("D43").Value = '0.6414'
$ws.Range("E43").Value = '  +3.13%  '

# Row 44
$ws.Range("D44").Value = '0.9959'
$ws.Range("E44").Value = '  -0.18%  '

# Row 45
$ws.Range("D45").Value = '13.82'
$ws.Range("E45").Value = '  +3.86%  '

# Row 46
$ws.Range("D46").Value = '0.6103'
$ws.Range("E46").Value = '  +1.91%  '

# Row 47
$ws.Range("D47").Value = '3.813'
$ws.Range("E47").Value = '  +0.09%  '

# Row 48
$ws.Range("D48").Value = '2.117'
$ws.Range("E48").Value = '  +3.58%  '

# Row 49
$ws.Range("D49").Value = '129.35'
$ws.Range("E49").Value = '  +1.22%  '

# Row 50
$ws.Range("D50").Value = '0.07270'
$ws.Range("E50").Value = '  +1.03%  '

# Row 51
$ws.Range("D51").Value = '79.13'
$ws.Range("E51").Value = '  +3.55%  '
